$wb = $excel.ActiveWorkbook

# --- Sheet2: "Schedule Number Counter" / "Data Recorder Index" updates ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = 135
$ws2.Range("B2").Value = 37

# --- Sheet3: append 9 new data rows (28-36) of "bad data" (unscheduled) ---
$ws3 = $wb.Worksheets.Item(3)

$newRows = @(
    @("3013696736", "1000004664", "13190563"),
    @("3013696737", "1000004665", "13190565"),
    @("3013696739", "1000004667", "13190569"),
    @("3013696742", "1000004669", "13190573"),
    @("3013696745", "1000004671", "13190576"),
    @("3013696746", "1000004672", "13190579"),
    @("3013696747", "1000004673", "13190581"),
    @("3013696748", "1000004674", "13190583"),
    @("3013696749", "1000004675", "13190585")
)

$r = 28
foreach ($row in $newRows) {
    $rangeA = $ws3.Range("A$r")
    $rangeA.NumberFormat = "@"
    $rangeA.Value = $row[0]
    $rangeA.Style = "Normal"

    $rangeB = $ws3.Range("B$r")
    $rangeB.NumberFormat = "@"
    $rangeB.Value = $row[1]
    $rangeB.Style = "Normal"

    $ws3.Range("C$r").Value = "schedNum"

    $rangeD = $ws3.Range("D$r")
    $rangeD.NumberFormat = "@"
    $rangeD.Value = $row[2]
    $rangeD.Style = "Normal"

    $r = $r + 1
}
